$wb = $excel.ActiveWorkbook

# --- Sheet "Project_Title" (4th sheet): new column of data, bold headers, wider columns, page setup ---
$ws4 = $wb.Worksheets.Item(4)

# Shift the existing "Engagement ..." value out to column C and put the new
# GE Healthcare entity name into column A (row 2), keeping column B as-is.
$ws4.Range("C2").Value = "Engagement GE Healthcare-FVA-101397"
$ws4.Range("A2").Value = "GE Healthcare-GE Healthcare Bio-Sciences AB-FVA-101397"

# Bold the header row (A1:B1)
$ws4.Range("A1:B1").Font.Bold = $true

# Resize the columns to fit their (now wider) content
$ws4.Columns.Item(1).ColumnWidth = 49.56
$ws4.Columns.Item(2).ColumnWidth = 15.78
$ws4.Columns.Item(3).ColumnWidth = 34.11

# Add print page setup (portrait) to the sheet
$ws4.PageSetup.Orientation = 1

# --- Update stored selections on every sheet (last one activated wins as the active tab) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C15").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D21").Select()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B6").Select()

# Project_Title becomes the active sheet/tab, selection on C8
$ws4.Activate()
$ws4.Range("C8").Select()
